# Scheduled-runner price refresh for the Leve profit-calculation sheets.
#
# Re-pulls current market-board averages (columns H: currentAveragePrice,
# I: currentAveragePriceNQ, J: currentAveragePriceHQ) and recomputes the
# dependent Leve price/profit columns (K: LevePriceNQ, L: LevePriceHQ,
# M: LeveProfitNQ, N: LeveProfitHQ) for the rows whose underlying item
# prices moved since the last run. Values are written as plain numbers
# (matching how the source sheet stores them -- no formulas), one row at
# a time per crafting-job worksheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 26000
$ws.Range("I47").Value = 25000
$ws.Range("J47").Value = 26500
$ws.Range("K47").Value = 25000
$ws.Range("L47").Value = 26500
$ws.Range("M47").Value = -24028
$ws.Range("N47").Value = -28444
$ws.Range("H107").Value = 1337
$ws.Range("I107").Value = 1707.7273
$ws.Range("J107").Value = 827.25
$ws.Range("K107").Value = 1707.7273
$ws.Range("L107").Value = 827.25
$ws.Range("M107").Value = 212.2727
$ws.Range("N107").Value = -4667.25
$ws.Range("H111").Value = 549.6667
$ws.Range("I111").Value = 549.6667
$ws.Range("K111").Value = 1649.0001
$ws.Range("M111").Value = 1417.9999
$ws.Range("H129").Value = 1271.1111
$ws.Range("I129").Value = 388
$ws.Range("J129").Value = 1610.7693
$ws.Range("K129").Value = 1164
$ws.Range("L129").Value = 4832.3079
$ws.Range("M129").Value = 3836
$ws.Range("N129").Value = -14832.3079
$ws.Range("H132").Value = 1990.9318
$ws.Range("I132").Value = 1670.1666
$ws.Range("J132").Value = 3434.375
$ws.Range("K132").Value = 5010.4998
$ws.Range("L132").Value = 10303.125
$ws.Range("M132").Value = -2480.4998
$ws.Range("N132").Value = -15363.125
$ws.Range("H137").Value = 2339.46
$ws.Range("I137").Value = 2155.111
$ws.Range("J137").Value = 2813.5
$ws.Range("K137").Value = 6465.333
$ws.Range("L137").Value = 8440.5
$ws.Range("M137").Value = -3915.333
$ws.Range("N137").Value = -13540.5
$ws.Range("H141").Value = 1268.8
$ws.Range("I141").Value = 860.58826
$ws.Range("J141").Value = 2136.25
$ws.Range("K141").Value = 2581.76478
$ws.Range("L141").Value = 6408.75
$ws.Range("M141").Value = 2598.23522
$ws.Range("N141").Value = -16768.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1160152.6
$ws.Range("I32").Value = 1245837
$ws.Range("J32").Value = 11982.8
$ws.Range("K32").Value = 1245837
$ws.Range("L32").Value = 11982.8
$ws.Range("M32").Value = -1245550
$ws.Range("N32").Value = -12556.8
$ws.Range("H63").Value = 3700
$ws.Range("I63").Value = 2950
$ws.Range("J63").Value = 4450
$ws.Range("K63").Value = 2950
$ws.Range("L63").Value = 4450
$ws.Range("M63").Value = -2264
$ws.Range("N63").Value = -5822
$ws.Range("H66").Value = 3700
$ws.Range("I66").Value = 2950
$ws.Range("J66").Value = 4450
$ws.Range("K66").Value = 14750
$ws.Range("L66").Value = 22250
$ws.Range("M66").Value = -11318
$ws.Range("N66").Value = -29114

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 25000
$ws.Range("I63").Value = 20000
$ws.Range("J63").Value = 30000
$ws.Range("K63").Value = 20000
$ws.Range("L63").Value = 30000
$ws.Range("M63").Value = -19314
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 25000
$ws.Range("I66").Value = 20000
$ws.Range("J66").Value = 30000
$ws.Range("K66").Value = 60000
$ws.Range("L66").Value = 90000
$ws.Range("M66").Value = -56568
$ws.Range("N66").Value = -96864
$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 920
$ws.Range("N107").Value = -5040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2126.0513
$ws.Range("I31").Value = 1521.18
$ws.Range("J31").Value = 3206.1785
$ws.Range("K31").Value = 1521.18
$ws.Range("L31").Value = 3206.1785
$ws.Range("M31").Value = -1226.18
$ws.Range("N31").Value = -3796.1785
$ws.Range("H34").Value = 2126.0513
$ws.Range("I34").Value = 1521.18
$ws.Range("J34").Value = 3206.1785
$ws.Range("K34").Value = 1521.18
$ws.Range("L34").Value = 3206.1785
$ws.Range("M34").Value = -1319.18
$ws.Range("N34").Value = -3610.1785
$ws.Range("H74").Value = 10941.777
$ws.Range("J74").Value = 10941.777
$ws.Range("L74").Value = 10941.777
$ws.Range("N74").Value = -12689.777
$ws.Range("H77").Value = 10941.777
$ws.Range("J77").Value = 10941.777
$ws.Range("L77").Value = 32825.331
$ws.Range("N77").Value = -41561.331
$ws.Range("H132").Value = 2024.9269
$ws.Range("I132").Value = 1073.9286
$ws.Range("K132").Value = 3221.7858
$ws.Range("M132").Value = -691.7857999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 198.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 198.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 595.5
$ws.Range("M7").ClearContents() | Out-Null
$ws.Range("N7").Value = -819.5
$ws.Range("H80").Value = 1283.5555
$ws.Range("J80").Value = 987.5
$ws.Range("L80").Value = 2962.5
$ws.Range("N80").Value = -4834.5
$ws.Range("H83").Value = 1283.5555
$ws.Range("J83").Value = 987.5
$ws.Range("L83").Value = 8887.5
$ws.Range("N83").Value = -18247.5
$ws.Range("H92").Value = 20833510
$ws.Range("I92").Value = 22727444
$ws.Range("J92").Value = 233.5
$ws.Range("K92").Value = 68182332
$ws.Range("L92").Value = 700.5
$ws.Range("M92").Value = -68181084
$ws.Range("N92").Value = -3196.5
$ws.Range("H97").Value = 519.125
$ws.Range("I97").Value = 650
$ws.Range("J97").Value = 475.5
$ws.Range("K97").Value = 1950
$ws.Range("L97").Value = 1426.5
$ws.Range("M97").Value = -1454
$ws.Range("N97").Value = -2418.5
$ws.Range("H124").Value = 3506.818
$ws.Range("I124").Value = 1326.3334
$ws.Range("J124").Value = 4324.5
$ws.Range("K124").Value = 3979.0002
$ws.Range("L124").Value = 12973.5
$ws.Range("M124").Value = 930.9998000000001
$ws.Range("N124").Value = -22793.5
$ws.Range("H131").Value = 1532.5167
$ws.Range("I131").Value = 1255.3636
$ws.Range("J131").Value = 1594.7347
$ws.Range("K131").Value = 3766.0908
$ws.Range("L131").Value = 4784.2041
$ws.Range("M131").Value = 1273.9092
$ws.Range("N131").Value = -14864.2041

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2338.9
$ws.Range("I113").Value = 2153.3333
$ws.Range("J113").Value = 2418.4285
$ws.Range("K113").Value = 2153.3333
$ws.Range("L113").Value = 2418.4285
$ws.Range("M113").Value = 16.66670000000022
$ws.Range("N113").Value = -6758.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2319.1667
$ws.Range("I82").Value = 1842.8334
$ws.Range("K82").Value = 1842.8334
$ws.Range("M82").Value = -1481.8334
$ws.Range("H85").Value = 2319.1667
$ws.Range("I85").Value = 1842.8334
$ws.Range("K85").Value = 1842.8334
$ws.Range("M85").Value = -594.8334
$ws.Range("H132").Value = 3930.4614
$ws.Range("I132").Value = 3589.6
$ws.Range("J132").Value = 4143.5
$ws.Range("K132").Value = 10768.8
$ws.Range("L132").Value = 12430.5
$ws.Range("M132").Value = -8238.799999999999
$ws.Range("N132").Value = -17490.5
$ws.Range("H137").Value = 21749.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 8000
$ws.Range("J39").Value = 7333.3335
$ws.Range("L39").Value = 7333.3335
$ws.Range("N39").Value = -8159.3335
$ws.Range("H113").Value = 594.0909
$ws.Range("I113").Value = 637.2222
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 1911.6666
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 258.3334
$ws.Range("N113").Value = -5540
$ws.Range("H126").Value = 1002.5
$ws.Range("I126").Value = 528.4666999999999
$ws.Range("K126").Value = 1585.4001
$ws.Range("M126").Value = 884.5999000000002
$ws.Range("H132").Value = 1731.7858
$ws.Range("I132").Value = 1134.0646
$ws.Range("J132").Value = 3416.2727
$ws.Range("K132").Value = 3402.1938
$ws.Range("L132").Value = 10248.8181
$ws.Range("M132").Value = -872.1938
$ws.Range("N132").Value = -15308.8181
$ws.Range("H136").Value = 24661098
$ws.Range("I136").Value = 35752104
$ws.Range("K136").Value = 107256312
$ws.Range("M136").Value = -107253762

